$wb = $excel.ActiveWorkbook

# --- Update the conversion summary text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.77 = 10365.98 pesos`n✅ 10365.98 pesos = 2.75 = 957.9 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 361.385
$wsTasas.Range("O10").Value = 3746.11
$wsTasas.Range("N12").Value = 3765.9
$wsTasas.Range("O12").Value = 348.001
